# Actualizacion automatica del mapa (2025-10-27 19:24:39)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update existing case status on row 6 (col G / "Estado")
$ws.Range("G6").Value = "Pendiente de Traspaso PROPIO"

# 2) Append a new case as row 75
#    Columns A-H, J-L, O-R are plain text in this sheet (even when the
#    text looks numeric, e.g. "Caso" = "-657"), so force text formatting
#    before writing, then reset the cell style back to Normal so no
#    extra style index is left behind on save.
$textRng1 = $ws.Range("A75:H75")
$textRng1.NumberFormat = "@"
$ws.Cells.Item(75, 1).Value  = "-657"
$ws.Cells.Item(75, 2).Value  = "10/27/2025"
$ws.Cells.Item(75, 3).Value  = "Conde 1632"
$ws.Cells.Item(75, 4).Value  = "13"
$ws.Cells.Item(75, 5).Value  = "810454540"
$ws.Cells.Item(75, 6).Value  = "NEW"
$ws.Cells.Item(75, 7).Value  = "Pendiente"
$ws.Cells.Item(75, 8).Value  = "Poste inclinado cambiar o desmontar"
$textRng1.Style = "Normal"

# Column I ("Attachments") is numeric
$ws.Cells.Item(75, 9).Value = 1

$textRng2 = $ws.Range("J75:L75")
$textRng2.NumberFormat = "@"
$ws.Cells.Item(75, 10).Value = "Cambio"
$ws.Cells.Item(75, 11).Value = "Sin equipos"
$ws.Cells.Item(75, 12).Value = "Poste"
$textRng2.Style = "Normal"

# Columns M/N ("Coordenada_X" / "Coordenada_Y") are numeric
$ws.Cells.Item(75, 13).Value = -58.461492
$ws.Cells.Item(75, 14).Value = -34.57199

$textRng3 = $ws.Range("O75:R75")
$textRng3.NumberFormat = "@"
$ws.Cells.Item(75, 15).Value = "Colegiales"
$ws.Cells.Item(75, 16).Value = "Capital Norte"
$ws.Cells.Item(75, 17).Value = "ATH-R"
$ws.Cells.Item(75, 18).Value = "Fuera de Poligono OVL"
$textRng3.Style = "Normal"
